# converter jin to kg unit
# Column N = "重量(公斤)" (weight in kg) on sheet "食材" (ingredients).
# The jin ("斤") quantities are converted to kilograms (1 jin = 0.6 kg) and
# written back as text (matching how the source data was originally stored).
# A leading apostrophe forces Excel to keep the numeric-looking text as a
# literal string instead of auto-converting it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value  = "'15.0"
$ws.Range("N3").Value  = "'5.3999999999999995"
$ws.Range("N4").Value  = "'1.7999999999999998"
$ws.Range("N5").Value  = "'1.2"
$ws.Range("N6").Value  = "'1.7999999999999998"
$ws.Range("N7").Value  = "'4.8"
$ws.Range("N8").Value  = "'4.8"
$ws.Range("N9").Value  = "'1.7999999999999998"
$ws.Range("N10").Value = "'1.2"
